# Apply the scripted updates to the worksheet.
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-10-27 Monday" "2025-10-28 Tuesday"

Replace-Text "14÷8=1, 6" "91÷8=11, 3"
Replace-Text "93÷3=31, 0" "40÷3=13, 1"
Replace-Text "25÷2=12, 1" "77÷4=19, 1"
Replace-Text "75÷8=9, 3" "15÷8=1, 7"
Replace-Text "89÷8=11, 1" "56÷9=6, 2"

Replace-Text "42÷7=6, 0" "22÷6=3, 4"
Replace-Text "47÷7=6, 5" "81÷8=10, 1"
Replace-Text "17÷4=4, 1" "93÷5=18, 3"
Replace-Text "38÷4=9, 2" "69÷3=23, 0"
Replace-Text "59÷8=7, 3" "25÷7=3, 4"

Replace-Text "54÷8=6, 6" "28÷7=4, 0"
Replace-Text "81÷6=13, 3" "22÷6=3, 4"
Replace-Text "13÷4=3, 1" "53÷9=5, 8"
Replace-Text "18÷7=2, 4" "33÷3=11, 0"
Replace-Text "15÷4=3, 3" "69÷8=8, 5"

Replace-Text "51÷6=8, 3" "57÷7=8, 1"
Replace-Text "98÷5=19, 3" "77÷7=11, 0"
Replace-Text "17÷3=5, 2" "27÷5=5, 2"
Replace-Text "40÷5=8, 0" "30÷9=3, 3"
Replace-Text "92÷7=13, 1" "30÷2=15, 0"

Replace-Text "24÷3=8, 0" "79÷7=11, 2"
Replace-Text "52÷2=26, 0" "23÷3=7, 2"
Replace-Text "63÷8=7, 7" "73÷9=8, 1"
Replace-Text "22÷8=2, 6" "91÷6=15, 1"
